$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) "Rak pilkoppling 8" (the first HTTP-arrow connector) loses its
#    arrow head at the start (headEnd -> none). msoArrowheadNone = 1.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Rak pilkoppling 8") {
        $sh.Line.BeginArrowheadStyle = 1
    }
}

# 2) Remove the three extra straight-arrow connectors
#    ("Rak pilkoppling 10/11/12") - they are replaced by a text label.
$namesToDelete = @("Rak pilkoppling 10", "Rak pilkoppling 11", "Rak pilkoppling 12")
foreach ($name in $namesToDelete) {
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -eq $name) {
            $sh.Delete()
        }
    }
}

# 3) Add a new textbox ("textruta 1") with the text "HTTP GET" where the
#    connectors used to be. PowerPoint's COM surface expresses shape
#    geometry in points (1 pt = 12700 EMU); the target OOXML positions
#    are given in EMU, so feed the point-equivalents in.
$tb = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$tb.Name = "textruta 1"
$tb.Left = 395.74952705905514
$tb.Top = 164.25850683700787
$tb.Width = 86.48141862283464
$tb.Height = 29.081259742519684

$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $false
$tb.TextFrame.AutoSize = 1
$tb.TextFrame.TextRange.Text = "HTTP GET"
